# Fix CF Memory Issue
#
# The snippet/warning counts in columns C (# of snippets with warnings) and
# D (# of warnings) were recounted, which shifts the dependent Kendall/Spearman
# correlation statistics in columns F-I (Tau, Kendall p-value, Rho, Spearman
# p-value). This updates the corrected values on the affected rows of the
# all_tools, checker_framework, and infer sheets (typestate_checker is
# unaffected).

$wb = $excel.ActiveWorkbook

$wsAll   = $wb.Worksheets.Item("all_tools")
$wsCf    = $wb.Worksheets.Item("checker_framework")
$wsInfer = $wb.Worksheets.Item("infer")


# --- all_tools ---
$wsAll.Cells.Item(2, 3).Value = 4   # C2
$wsAll.Cells.Item(2, 4).Value = 21   # D2
$wsAll.Cells.Item(2, 6).Value = 0.09719878783949158   # F2
$wsAll.Cells.Item(2, 7).Value = 0.5729099316418407   # G2
$wsAll.Cells.Item(2, 8).Value = 0.1360629871783227   # H2
$wsAll.Cells.Item(2, 9).Value = 0.5358931384723264   # I2
$wsAll.Cells.Item(3, 3).Value = 4   # C3
$wsAll.Cells.Item(3, 4).Value = 21   # D3
$wsAll.Cells.Item(3, 6).Value = -0.09266821663323266   # F3
$wsAll.Cells.Item(3, 7).Value = 0.5984721667632871   # G3
$wsAll.Cells.Item(3, 8).Value = -0.1263163607514111   # H3
$wsAll.Cells.Item(3, 9).Value = 0.5657528992105718   # I3
$wsAll.Cells.Item(4, 3).Value = 4   # C4
$wsAll.Cells.Item(4, 4).Value = 21   # D4
$wsAll.Cells.Item(4, 6).Value = -0.1259696449729016   # F4
$wsAll.Cells.Item(4, 7).Value = 0.4679676774681618   # G4
$wsAll.Cells.Item(4, 8).Value = -0.1721605113945802   # H4
$wsAll.Cells.Item(4, 9).Value = 0.4321598963432832   # I4
$wsAll.Cells.Item(9, 4).Value = 395   # D9
$wsAll.Cells.Item(9, 6).Value = -0.2278622001693345   # F9
$wsAll.Cells.Item(9, 7).Value = 0.001503036931867509   # G9
$wsAll.Cells.Item(9, 8).Value = -0.3234568550885335   # H9
$wsAll.Cells.Item(9, 9).Value = 0.0010284898212685   # I9
$wsAll.Cells.Item(13, 4).Value = 30   # D13
$wsAll.Cells.Item(13, 6).Value = -0.1130105839368017   # F13
$wsAll.Cells.Item(13, 7).Value = 0.4404006981390032   # G13
$wsAll.Cells.Item(13, 8).Value = -0.1377826417184593   # H13
$wsAll.Cells.Item(13, 9).Value = 0.4677936200616143   # I13
$wsAll.Cells.Item(14, 4).Value = 30   # D14
$wsAll.Cells.Item(14, 6).Value = -0.1668825406867738   # F14
$wsAll.Cells.Item(14, 7).Value = 0.2550672974746294   # G14
$wsAll.Cells.Item(14, 8).Value = -0.2023907691537807   # H14
$wsAll.Cells.Item(14, 9).Value = 0.2834571786891915   # I14
$wsAll.Cells.Item(15, 4).Value = 30   # D15
$wsAll.Cells.Item(15, 6).Value = 0.5876550364713691   # F15
$wsAll.Cells.Item(15, 7).Value = 0.00006023530742288514   # G15
$wsAll.Cells.Item(15, 8).Value = 0.7491931143441224   # H15
$wsAll.Cells.Item(15, 9).Value = 0.000001904021591222836   # I15
$wsAll.Cells.Item(16, 4).Value = 30   # D16
$wsAll.Cells.Item(16, 6).Value = 0.5876550364713691   # F16
$wsAll.Cells.Item(16, 7).Value = 0.00006023530742288514   # G16
$wsAll.Cells.Item(16, 8).Value = 0.7491931143441224   # H16
$wsAll.Cells.Item(16, 9).Value = 0.000001904021591222836   # I16

# --- checker_framework ---
$wsCf.Cells.Item(2, 3).Value = 3   # C2
$wsCf.Cells.Item(2, 4).Value = 7   # D2
$wsCf.Cells.Item(2, 6).Value = 0.2297034206521828   # F2
$wsCf.Cells.Item(2, 7).Value = 0.187290841633909   # G2
$wsCf.Cells.Item(2, 8).Value = 0.2809695424230301   # H2
$wsCf.Cells.Item(2, 9).Value = 0.1940516726155841   # I2
$wsCf.Cells.Item(3, 3).Value = 3   # C3
$wsCf.Cells.Item(3, 4).Value = 7   # D3
$wsCf.Cells.Item(3, 6).Value = -0.2277100170213244   # F3
$wsCf.Cells.Item(3, 7).Value = 0.2003280221869526   # G3
$wsCf.Cells.Item(3, 8).Value = -0.2719723502938716   # H3
$wsCf.Cells.Item(3, 9).Value = 0.2093250956596323   # I3
$wsCf.Cells.Item(4, 3).Value = 3   # C4
$wsCf.Cells.Item(4, 4).Value = 7   # D4
$wsCf.Cells.Item(4, 6).Value = -0.263477777620917   # F4
$wsCf.Cells.Item(4, 7).Value = 0.1329850671160174   # G4
$wsCf.Cells.Item(4, 8).Value = -0.3201778730528596   # H4
$wsCf.Cells.Item(4, 9).Value = 0.1363764100850406   # I4
$wsCf.Cells.Item(9, 3).Value = 19   # C9
$wsCf.Cells.Item(9, 4).Value = 52   # D9
$wsCf.Cells.Item(9, 6).Value = -0.2289041597670328   # F9
$wsCf.Cells.Item(9, 7).Value = 0.004154062204876697   # G9
$wsCf.Cells.Item(9, 8).Value = -0.2872671746617843   # H9
$wsCf.Cells.Item(9, 9).Value = 0.003756720549751365   # I9
$wsCf.Cells.Item(13, 3).Value = 12   # C13
$wsCf.Cells.Item(13, 4).Value = 12   # D13
$wsCf.Cells.Item(13, 6).Value = -0.06524669105718928   # F13
$wsCf.Cells.Item(13, 7).Value = 0.6720517686591893   # G13
$wsCf.Cells.Item(13, 8).Value = -0.07861110510609637   # H13
$wsCf.Cells.Item(13, 9).Value = 0.679668768010645   # I13
$wsCf.Cells.Item(14, 3).Value = 12   # C14
$wsCf.Cells.Item(14, 4).Value = 12   # D14
$wsCf.Cells.Item(14, 6).Value = -0.153506269493634   # F14
$wsCf.Cells.Item(14, 7).Value = 0.3197630116677985   # G14
$wsCf.Cells.Item(14, 8).Value = -0.1847566494939108   # H14
$wsCf.Cells.Item(14, 9).Value = 0.328369658424491   # I14
$wsCf.Cells.Item(15, 3).Value = 12   # C15
$wsCf.Cells.Item(15, 4).Value = 12   # D15
$wsCf.Cells.Item(15, 6).Value = 0.443677499188887   # F15
$wsCf.Cells.Item(15, 7).Value = 0.003993565524162987   # G15
$wsCf.Cells.Item(15, 8).Value = 0.5345555147214552   # H15
$wsCf.Cells.Item(15, 9).Value = 0.002341385142708154   # I15
$wsCf.Cells.Item(16, 3).Value = 12   # C16
$wsCf.Cells.Item(16, 4).Value = 12   # D16
$wsCf.Cells.Item(16, 6).Value = 0.443677499188887   # F16
$wsCf.Cells.Item(16, 7).Value = 0.003993565524162987   # G16
$wsCf.Cells.Item(16, 8).Value = 0.5345555147214552   # H16
$wsCf.Cells.Item(16, 9).Value = 0.002341385142708154   # I16

# --- infer ---
$wsInfer.Cells.Item(13, 3).Value = 18   # C13
$wsInfer.Cells.Item(13, 6).Value = -0.1435427203258164   # F13
$wsInfer.Cells.Item(13, 7).Value = 0.3516806827985527   # G13
$wsInfer.Cells.Item(13, 8).Value = -0.172944431233412   # H13
$wsInfer.Cells.Item(13, 9).Value = 0.360759748541673   # I13
$wsInfer.Cells.Item(14, 3).Value = 18   # C14
$wsInfer.Cells.Item(14, 6).Value = -0.153506269493634   # F14
$wsInfer.Cells.Item(14, 7).Value = 0.3197630116677985   # G14
$wsInfer.Cells.Item(14, 8).Value = -0.1847566494939108   # H14
$wsInfer.Cells.Item(14, 9).Value = 0.328369658424491   # I14
$wsInfer.Cells.Item(15, 3).Value = 18   # C15
$wsInfer.Cells.Item(15, 6).Value = 0.6916149252062063   # F15
$wsInfer.Cells.Item(15, 7).Value = 0.000007211867226751588   # G15
$wsInfer.Cells.Item(15, 8).Value = 0.8332777141246215   # H15
$wsInfer.Cells.Item(15, 9).Value = 0.00000001098076308139691   # I15
$wsInfer.Cells.Item(16, 3).Value = 18   # C16
$wsInfer.Cells.Item(16, 6).Value = 0.6916149252062063   # F16
$wsInfer.Cells.Item(16, 7).Value = 0.000007211867226751588   # G16
$wsInfer.Cells.Item(16, 8).Value = 0.8332777141246215   # H16
$wsInfer.Cells.Item(16, 9).Value = 0.00000001098076308139691   # I16
